$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Forecast Comparison sheet - Seasonality Index (L) and Inventory Coverage (H) updates
$ws1.Range("L2").Value = 0.86
$ws1.Range("L3").Value = 1.17
$ws1.Range("L4").Value = 0.8100000000000001
$ws1.Range("L5").Value = 1.15

$ws1.Range("H6").Value = 44
$ws1.Range("L6").Value = 0.89

$ws1.Range("H7").Value = 43
$ws1.Range("L7").Value = 1.1

$ws1.Range("H8").Value = 28
$ws1.Range("L8").Value = 1.07

$ws1.Range("H9").Value = 27
$ws1.Range("L9").Value = 0.86

$ws1.Range("H10").Value = 26
$ws1.Range("L10").Value = 0.85

$ws1.Range("H11").Value = 37.5
$ws1.Range("L11").Value = 1.15

$ws1.Range("H12").Value = 24.33
$ws1.Range("L12").Value = 0.85

$ws1.Range("H13").Value = 23.33
$ws1.Range("L13").Value = 1.07

$ws1.Range("H14").Value = 22.33
$ws1.Range("L14").Value = 1.15

$ws1.Range("H15").Value = 21.33
$ws1.Range("L15").Value = 0.91

$ws1.Range("H16").Value = 20.33
$ws1.Range("L16").Value = 1.18

$ws1.Range("H17").Value = 14.5
$ws1.Range("L17").Value = 0.91

# Summary sheet updates (keep these as text values, matching existing column format)
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "4"
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "0"
